$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44594
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 120
$ws.Range("O2").Value = 2800
$ws.Range("P2").Value = 2650
$ws.Range("R2").Value = 'Provincia de Linares'
$ws.Range("S2").Value = 1325
# Row 3
$ws.Range("D3").Value = 44970
$ws.Range("M3").Value = 50
# Row 4
$ws.Range("D4").Value = 44970
$ws.Range("L4").Value = 'Segunda'
$ws.Range("N4").Value = 2500
$ws.Range("O4").Value = 2500
$ws.Range("P4").Value = 2500
$ws.Range("S4").Value = 1250
# Row 5
$ws.Range("D5").Value = 44596
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 120
$ws.Range("O5").Value = 2700
$ws.Range("P5").Value = 2600
$ws.Range("R5").Value = 'Provincia de Linares'
$ws.Range("S5").Value = 1300
# Row 6
$ws.Range("D6").Value = 44951
$ws.Range("N6").Value = 2800
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = 2900
$ws.Range("S6").Value = 1450
# Row 8
$ws.Range("D8").Value = 44960
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 2500
$ws.Range("O8").Value = 2500
$ws.Range("P8").Value = 2500
$ws.Range("S8").Value = 1250
# Row 9
$ws.Range("D9").Value = 44988
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("R9").Value = 'Provincia de Diguillín'
$ws.Range("S9").Value = 1500
# Row 10
$ws.Range("D10").Value = 44988
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("Q10").Value = '$/bandeja 2 kilos'
$ws.Range("S10").Value = 1250
$ws.Range("T10").Value = 2
# Row 11
$ws.Range("D11").Value = 45009
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 4000
$ws.Range("O11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("R11").Value = 'Provincia de Diguillín'
$ws.Range("S11").Value = 2000
# Row 12
$ws.Range("D12").Value = 44974
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 3000
$ws.Range("P12").Value = 3000
$ws.Range("S12").Value = 1500
# Row 13
$ws.Range("D13").Value = 44974
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 30
# Row 14
$ws.Range("D14").Value = 44979
$ws.Range("M14").Value = 30
# Row 15
$ws.Range("D15").Value = 44979
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 2500
$ws.Range("O15").Value = 2500
$ws.Range("P15").Value = 2500
$ws.Range("R15").Value = 'Provincia de Diguillín'
$ws.Range("S15").Value = 1250
# Row 16
$ws.Range("D16").Value = 44949
$ws.Range("N16").Value = 2800
$ws.Range("O16").Value = 3000
$ws.Range("P16").Value = 2900
$ws.Range("S16").Value = 1450
# Row 17
$ws.Range("D17").Value = 44985
# Row 18
$ws.Range("D18").Value = 44985
$ws.Range("M18").Value = 50
# Row 19
$ws.Range("D19").Value = 44174
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 3700
$ws.Range("O19").Value = 3800
$ws.Range("P19").Value = 3747
$ws.Range("S19").Value = 1874
# Row 20
$ws.Range("D20").Value = 44935
# Row 21
$ws.Range("D21").Value = 44944
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 2500
$ws.Range("O21").Value = 2500
$ws.Range("P21").Value = 2500
$ws.Range("S21").Value = 1250
# Row 22
$ws.Range("D22").Value = 44956
# Row 23
$ws.Range("D23").Value = 44952
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 30
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("S23").Value = 1500
# Row 24
$ws.Range("D24").Value = 44965
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 3000
$ws.Range("O24").Value = 3000
$ws.Range("P24").Value = 3000
$ws.Range("S24").Value = 1500
# Row 25
$ws.Range("D25").Value = 45008
$ws.Range("M25").Value = 30
$ws.Range("N25").Value = 4000
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 4000
$ws.Range("R25").Value = 'Provincia de Diguillín'
$ws.Range("S25").Value = 2000
# Row 26
$ws.Range("D26").Value = 44963
$ws.Range("M26").Value = 50
# Row 27
$ws.Range("D27").Value = 44963
$ws.Range("L27").Value = 'Segunda'
$ws.Range("N27").Value = 2500
$ws.Range("O27").Value = 2500
$ws.Range("P27").Value = 2500
$ws.Range("S27").Value = 1250
# Row 28
$ws.Range("D28").Value = 44992
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 3000
$ws.Range("O28").Value = 3000
$ws.Range("P28").Value = 3000
$ws.Range("S28").Value = 1500
# Row 29
$ws.Range("D29").Value = 44966
$ws.Range("L29").Value = 'Segunda'
$ws.Range("N29").Value = 2500
$ws.Range("O29").Value = 2500
$ws.Range("P29").Value = 2500
$ws.Range("S29").Value = 1250
# Row 30
$ws.Range("D30").Value = 45006
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 4000
$ws.Range("O30").Value = 4000
$ws.Range("P30").Value = 4000
$ws.Range("S30").Value = 2000
# Row 31
$ws.Range("D31").Value = 44942
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 2500
$ws.Range("O31").Value = 2500
$ws.Range("P31").Value = 2500
$ws.Range("S31").Value = 1250
# Row 32
$ws.Range("D32").Value = 44953
$ws.Range("L32").Value = 'Primera'
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 3000
$ws.Range("P32").Value = 3000
$ws.Range("S32").Value = 1500
# Row 33
$ws.Range("D33").Value = 44932
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 3000
$ws.Range("P33").Value = 3000
$ws.Range("S33").Value = 1500
# Row 34
$ws.Range("D34").Value = 44967
$ws.Range("M34").Value = 50
# Row 35
$ws.Range("D35").Value = 44967
$ws.Range("L35").Value = 'Segunda'
$ws.Range("N35").Value = 2500
$ws.Range("O35").Value = 2500
$ws.Range("P35").Value = 2500
$ws.Range("S35").Value = 1250
# Row 36
$ws.Range("D36").Value = 44971
# Row 37
$ws.Range("D37").Value = 44972
# Row 38
$ws.Range("D38").Value = 44187
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 2800
$ws.Range("P38").Value = 2900
$ws.Range("R38").Value = 'Provincia de Linares'
$ws.Range("S38").Value = 1450
# Row 39
$ws.Range("D39").Value = 44187
$ws.Range("L39").Value = 'Primera'
$ws.Range("M39").Value = 65
$ws.Range("N39").Value = 1400
$ws.Range("O39").Value = 1500
$ws.Range("P39").Value = 1446
$ws.Range("Q39").Value = '$/envase 1 kilo'
$ws.Range("S39").Value = 1446
$ws.Range("T39").Value = 1
# Row 40
$ws.Range("D40").Value = 44937
$ws.Range("M40").Value = 100
$ws.Range("N40").Value = 2500
$ws.Range("O40").Value = 3000
$ws.Range("P40").Value = 2750
$ws.Range("S40").Value = 1375
# Row 41
$ws.Range("D41").Value = 44994
$ws.Range("M41").Value = 60
$ws.Range("N41").Value = 3000
$ws.Range("O41").Value = 3200
$ws.Range("P41").Value = 3100
$ws.Range("Q41").Value = '$/bandeja 2 kilos'
$ws.Range("S41").Value = 1550
$ws.Range("T41").Value = 2
# Row 42
$ws.Range("D42").Value = 44931
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 3000
$ws.Range("P42").Value = 3000
$ws.Range("S42").Value = 1500
# Row 43
$ws.Range("D43").Value = 44540
$ws.Range("M43").Value = 240
$ws.Range("N43").Value = 3500
$ws.Range("O43").Value = 3800
$ws.Range("P43").Value = 3650
$ws.Range("R43").Value = 'Región del Maule'
$ws.Range("S43").Value = 1825
# Row 44
$ws.Range("D44").Value = 44181
$ws.Range("L44").Value = 'Primera'
$ws.Range("M44").Value = 65
$ws.Range("N44").Value = 3600
$ws.Range("O44").Value = 3800
$ws.Range("P44").Value = 3692
$ws.Range("S44").Value = 1846
# Row 45
$ws.Range("D45").Value = 44181
$ws.Range("M45").Value = 80
$ws.Range("N45").Value = 1800
$ws.Range("O45").Value = 2000
$ws.Range("P45").Value = 1875
$ws.Range("Q45").Value = '$/envase 1 kilo'
$ws.Range("S45").Value = 1875
$ws.Range("T45").Value = 1
